$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current D:T values for every data row before rewriting,
# since the update reorders (permutes) the rows.
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $snapshot[$r] = $ws.Range("D$r`:T$r").Value2
}

# Write each destination row from the snapshot of its source row.
$ws.Range("D2:T2").Value2 = $snapshot[7]
$ws.Range("D3:T3").Value2 = $snapshot[4]
$ws.Range("D4:T4").Value2 = $snapshot[12]
$ws.Range("D5:T5").Value2 = $snapshot[2]
$ws.Range("D6:T6").Value2 = $snapshot[17]
$ws.Range("D7:T7").Value2 = $snapshot[15]
$ws.Range("D8:T8").Value2 = $snapshot[16]
$ws.Range("D9:T9").Value2 = $snapshot[22]
$ws.Range("D10:T10").Value2 = $snapshot[3]
$ws.Range("D11:T11").Value2 = $snapshot[23]
$ws.Range("D12:T12").Value2 = $snapshot[24]
$ws.Range("D13:T13").Value2 = $snapshot[19]
$ws.Range("D14:T14").Value2 = $snapshot[13]
$ws.Range("D15:T15").Value2 = $snapshot[25]
$ws.Range("D16:T16").Value2 = $snapshot[8]
$ws.Range("D17:T17").Value2 = $snapshot[9]
$ws.Range("D18:T18").Value2 = $snapshot[11]
$ws.Range("D19:T19").Value2 = $snapshot[6]
$ws.Range("D20:T20").Value2 = $snapshot[5]
$ws.Range("D21:T21").Value2 = $snapshot[18]
$ws.Range("D22:T22").Value2 = $snapshot[20]
$ws.Range("D23:T23").Value2 = $snapshot[21]
$ws.Range("D24:T24").Value2 = $snapshot[14]
$ws.Range("D25:T25").Value2 = $snapshot[10]
